$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three run timestamps in column B (Date)
$ws.Range("B2").Value = "Thu Dec 07 21:42:19 EST 2023"
$ws.Range("B3").Value = "Thu Dec 07 21:42:34 EST 2023"
$ws.Range("B5").Value = "Thu Dec 07 21:42:47 EST 2023"

# Mark the Extension Payments row (row 4) as not to be executed
$ws.Range("C4").Value = "DONOTRUN"

# Update the active selection to C4, matching the final saved state
$ws.Range("C4").Select()

# Widen column C to fit the new "DONOTRUN" text and disable bestFit/autofit sizing
$ws.Columns("C").ColumnWidth = 14.6667
